# Thesis progress workbook update — "Finished graphical part of background!"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core data edit: the timeline's pages-per-day target changed ---
# L16 (Start) moves from 46 to 53 actual pages written by that date,
# and N16 (Tdays) formula changes from 24-7 to 22-7 (fewer remaining days).
$ws.Range("L16").Value = 53
$ws.Range("N16").Formula = "=22-7"

# --- Row 22: the stray " " placeholder in E22 is removed so F22/G22 ---
# --- stop erroring out (#VALUE!) and compute real numbers again.    ---
$ws.Range("E22").ClearContents()

# --- Background section labels (I column) ---
# New marker on row 32 for the newly finished "Evaluation" section
$ws.Range("I32").Value = "Evaluation"
# Row 25 marker relabelled from the old lowercase "evaluation" to "Design"
$ws.Range("I25").Value = "Design"

# --- New J column: a secondary day-count series running alongside the ---
# --- existing schedule (J22=55, then +2 each subsequent day).        ---
$ws.Range("J22").Value = 55
$ws.Range("J23").Formula = "=J22+2"
$ws.Range("J24:J37").Formula = "=J23+2"

# --- View state: selection moved from E22 to C21 ---
$ws.Range("C21").Select() | Out-Null
